# Raw and Clean Data from SSA for July 27th
# Appends a new row (58) to the "out_vars" sheet with the date label
# "2020-07-27" (as text, matching the existing column-A convention) and
# the corresponding Confirmados/Negativos/Sospechosos/Defunciones counts
# plus the hospitalized percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 58

# Write the date as a text formula result, then flatten it to a static
# value via copy/paste-values. This avoids Excel's automatic text->date
# conversion (which would otherwise turn "2020-07-27" into a date serial)
# while keeping the cell a plain text value (no extra number-format style).
$ws.Cells.Item($newRow, 1).Formula = '="2020-07-27"'
$ws.Cells.Item($newRow, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4163)

$ws.Cells.Item($newRow, 2).Value = 395489
$ws.Cells.Item($newRow, 3).Value = 442884
$ws.Cells.Item($newRow, 4).Value = 85986
$ws.Cells.Item($newRow, 5).Value = 44022
$ws.Cells.Item($newRow, 6).Value = 27.72
